$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Brócoli" (Macroferia Regional de Talca) is
# inserted at row 517, pushing every existing record from row 517 onward
# down by one row (old row 517 -> new row 518, ..., old row 595 -> new row 596).

$ws.Rows(517).Insert()

# Seed the new row with the same recurring fields as the record right below it
# (post-shift row 518 = what used to be row 517), then override the fields
# that differ for this new entry (date + prices).
$ws.Range("A518:R518").Copy()
$ws.Range("A517").PasteSpecial()

$ws.Range("D517").Value2 = 45131
$ws.Range("K517").Value2 = 600
$ws.Range("L517").Value2 = 600
$ws.Range("M517").Value2 = 600
$ws.Range("P517").Value2 = 600
